$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 7711
$ws.Range("F3").Value = 3623
$ws.Range("F4").Value = 3924
$ws.Range("F6").Value = 107
$ws.Range("F8").Value = 130
$ws.Range("F9").Value = 203
$ws.Range("F10").Value = 538
$ws.Range("F11").Value = 22
$ws.Range("F12").Value = 181
$ws.Range("F13").Value = 17
$ws.Range("F15").Value = 23
$ws.Range("F17").Value = 368
$ws.Range("F18").Value = 4318
$ws.Range("F19").Value = 4318
$ws.Range("F21").Value = 431
$ws.Range("F22").Value = 1047
$ws.Range("F23").Value = 551
$ws.Range("F24").Value = 3480
$ws.Range("F26").Value = 119
$ws.Range("F27").Value = 3159
$ws.Range("F28").Value = 2465
$ws.Range("F29").Value = 85
$ws.Range("F30").Value = 91
$ws.Range("F31").Value = 3
$ws.Range("F33").Value = 141
$ws.Range("F35").Value = 53
$ws.Range("F36").Value = 44
$ws.Range("F37").Value = 117
$ws.Range("F38").Value = 4639
$ws.Range("F39").Value = 579
$ws.Range("F40").Value = 351
$ws.Range("F41").Value = 65
$ws.Range("F43").Value = 909
$ws.Range("F44").Value = 300
$ws.Range("F46").Value = 1752
$ws.Range("F47").Value = 275
$ws.Range("F48").Value = 46
$ws.Range("F49").Value = 634
$ws.Range("F50").Value = 757

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 5
$ws.Range("F4").Value = 453
$ws.Range("F6").Value = 5
$ws.Range("F8").Value = 70
$ws.Range("F9").Value = 109
$ws.Range("F17").Value = 111
$ws.Range("F24").Value = 663

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 7711
$ws.Range("F4").Value = 3623
$ws.Range("F5").Value = 3924
$ws.Range("F7").Value = 107
$ws.Range("F9").Value = 130
$ws.Range("F10").Value = 203
$ws.Range("F11").Value = 5
$ws.Range("F12").Value = 538
$ws.Range("F14").Value = 181
$ws.Range("F16").Value = 23
$ws.Range("F17").Value = 368
$ws.Range("F18").Value = 4318
$ws.Range("F19").Value = 4318
$ws.Range("F23").Value = 431
$ws.Range("F24").Value = 1047
$ws.Range("F25").Value = 551
$ws.Range("F26").Value = 3482
$ws.Range("F28").Value = 119
$ws.Range("F29").Value = 3159
$ws.Range("F30").Value = 2465
$ws.Range("F31").Value = 85
$ws.Range("F32").Value = 91
$ws.Range("F33").Value = 141
$ws.Range("F36").Value = 44
$ws.Range("F37").Value = 117
$ws.Range("F39").Value = 4639
$ws.Range("F41").Value = 579
$ws.Range("F42").Value = 351
$ws.Range("F44").Value = 65
$ws.Range("F45").Value = 909
$ws.Range("F46").Value = 300
$ws.Range("F47").Value = 1752
$ws.Range("F48").Value = 275
$ws.Range("F49").Value = 634
$ws.Range("F50").Value = 757
